# Aufwand_git.xlsx - add two new tracked work entries (2024-03-06) and
# update the view so the newly added rows are visible/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# New row 51: continuation of the "Bug Blank Nodes Anzeige" task (3h)
$ws.Range("A51").Value = 45357
$ws.Range("A51").NumberFormat = $ws.Range("A50").NumberFormat
$ws.Range("B51").Value = 3
$ws.Range("C51").Value = "Bug Blank Nodes Anzeige"

# New row 52: "Filter Overview" work item (5h) - new unique string
$ws.Range("A52").Value = 45357
$ws.Range("A52").NumberFormat = $ws.Range("A50").NumberFormat
$ws.Range("B52").Value = 5
$ws.Range("C52").Value = "Filter Overview"

# Move the view so the freshly entered rows are in frame and selected,
# matching the author's on-screen state after the edit.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C54").Select()
